# Fernandez21_data.xlsx - "updated datasets for rev 2"
#
# The data-column layout shifts one slot to the left starting at G:
#   - G (old "experiment id" text, e.g. "pw"/"freq"/"dur") becomes a
#     constant 1
#   - H picks up the old I value
#   - I picks up the old J value
#   - J becomes the ratio formula 10^-6*I/E (previously this formula
#     lived in K and referenced J)
#   - K (the old formula column) is cleared out entirely
# The header row (row 1) is relabelled to match, and a new "sem" header
# lands on K1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) -------------------------------------------------
$ws.Range("G1").Value = "ppc"
$ws.Range("H1").Value = "dur"
$ws.Range("I1").Value = "charge"
$ws.Range("J1").Value = "amp"
$ws.Range("K1").Value = "sem"

# ---- Data rows 2-11: capture old H/I/J values before overwriting --------
$oldI = @{}
$oldJ = @{}
for ($r = 2; $r -le 11; $r++) {
    $oldI[$r] = $ws.Cells.Item($r, 9).Value()
    $oldJ[$r] = $ws.Cells.Item($r, 10).Value()
}

# G becomes a constant 1 for every data row
$ws.Range("G2:G11").Value = 1

# H <- old I, I <- old J
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = $oldI[$r]
    $ws.Cells.Item($r, 9).Value = $oldJ[$r]
}

# J becomes the ratio formula. Row 2 stays a standalone formula; rows
# 3-11 are entered as one fill so the engine compresses them into a
# shared formula, mirroring how the original K3:K11 formula was shared.
$ws.Range("J2").Formula = "=10^-6*I2/E2"
$ws.Range("J3:J11").Formula = "=10^-6*I3/E3"

# K is cleared entirely (content removed from every data row)
$ws.Range("K2:K11").ClearContents()

# ---- Column width: the old bestFit width on K now belongs to J ----------
$ws.Columns.Item(10).ColumnWidth = 10.85

# ---- Selection moved from D2:D11 to D3:D11 -------------------------------
$ws.Range("D3:D11").Select()
